$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.310.62'
$ws.Range("D3").Value = '''1.874.99'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '''241.77'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.07883'
$ws.Range("E8").Value = '  +2.15%  '
$ws.Range("D9").Value = '''0.3122'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").Value = '''25.23'
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("D11").Value = '''0.08393'
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '''1.878.33'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '''5.243'
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = '''0.7175'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = '''91.26'
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").Value = '''6.204'
$ws.Range("E16").Value = '  +3.75%  '
$ws.Range("D17").Value = '''0.000008364'
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '''29.310.64'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '''240.68'
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("D20").Value = '''13.23'
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").Value = '''2.123.04'
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("D22").Value = '''0.9998'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '''7.800'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '''0.1593'
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").Value = '''162.67'
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").Value = '''9.050'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = '''18.53'
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").Value = '''1.506'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '''4.424'
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("D31").Value = '''4.344'
$ws.Range("D32").Value = '''1.186'
$ws.Range("E32").Value = '  -8.31%  '
$ws.Range("D33").Value = '''0.05355'
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("D34").Value = '''1.944'
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("D36").Value = '''1.177'
$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("D37").Value = '''2.697'
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("D38").Value = '''1.290.74'
$ws.Range("E38").Value = '  +11.62%  '
$ws.Range("D39").Value = '''0.01884'
$ws.Range("E39").Value = '  +1.45%  '
$ws.Range("D40").Value = '''2.735'
$ws.Range("E40").Value = '  +0.64%  '
$ws.Range("D41").Value = '''6.582'
$ws.Range("E41").Value = '  +3.43%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.8967'
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''110.65'
$ws.Range("E43").Value = '  +4.43%  '
$ws.Range("D44").Value = '''73.22'
$ws.Range("E45").Value = '  +8.91%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").Value = '''2.021.73'
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").Value = '''1.801'
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("D49").Value = '''0.5200'
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").Value = '''9.454'
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("D51").Value = '''0.4357'
$ws.Range("E51").Value = '  +1.29%  '
